$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 154.32
$ws.Range("I15").Value = 154.32
$ws.Range("K15").Value = 462.96
$ws.Range("M15").Value = -293.96

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 470
$ws.Range("I101").Value = 430.8
$ws.Range("J101").Value = 666
$ws.Range("K101").Value = 1292.4
$ws.Range("L101").Value = 1998
$ws.Range("M101").Value = 329.5999999999999
$ws.Range("N101").Value = -5242

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2175.5
$ws.Range("J112").Value = 2644.375
$ws.Range("L112").Value = 7933.125
$ws.Range("N112").Value = -10149.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1744.5
$ws.Range("I116").Value = 1659.5555
$ws.Range("J116").Value = 1999.3334
$ws.Range("K116").Value = 1659.5555
$ws.Range("L116").Value = 1999.3334
$ws.Range("M116").Value = 1782.4445
$ws.Range("N116").Value = -8883.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 4017
$ws.Range("J129").Value = 969.96295
$ws.Range("L129").Value = 2909.88885
$ws.Range("N129").Value = -12909.88885

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5005804
$ws.Range("I132").Value = 5958647
$ws.Range("J132").Value = 3376.5
$ws.Range("K132").Value = 17875941
$ws.Range("L132").Value = 10129.5
$ws.Range("M132").Value = -17873411
$ws.Range("N132").Value = -15189.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4589.7446
$ws.Range("I138").Value = 2773.25
$ws.Range("J138").Value = 5527.2905
$ws.Range("K138").Value = 8319.75
$ws.Range("L138").Value = 16581.8715
$ws.Range("M138").Value = -3179.75
$ws.Range("N138").Value = -26861.8715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 39124.367
$ws.Range("I32").Value = 17585.594
$ws.Range("J32").Value = 113876.586
$ws.Range("K32").Value = 17585.594
$ws.Range("L32").Value = 113876.586
$ws.Range("M32").Value = -17298.594
$ws.Range("N32").Value = -114450.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 10600
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10600
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 10600
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -11572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2565.7058
$ws.Range("I61").Value = 2527.25
$ws.Range("J61").Value = 2658
$ws.Range("K61").Value = 2527.25
$ws.Range("L61").Value = 2658
$ws.Range("M61").Value = -2315.25
$ws.Range("N61").Value = -3082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1677.8
$ws.Range("I74").Value = 1598.9231
$ws.Range("J74").Value = 1824.2858
$ws.Range("K74").Value = 1598.9231
$ws.Range("L74").Value = 1824.2858
$ws.Range("M74").Value = -724.9231
$ws.Range("N74").Value = -3572.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1677.8
$ws.Range("I77").Value = 1598.9231
$ws.Range("J77").Value = 1824.2858
$ws.Range("K77").Value = 7994.6155
$ws.Range("L77").Value = 9121.429
$ws.Range("M77").Value = -3626.6155
$ws.Range("N77").Value = -17857.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2565.7058
$ws.Range("I136").Value = 2527.25
$ws.Range("J136").Value = 2658
$ws.Range("K136").Value = 7581.75
$ws.Range("L136").Value = 7974
$ws.Range("M136").Value = -5031.75
$ws.Range("N136").Value = -13074

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 425.7143
$ws.Range("I22").Value = 463.33334
$ws.Range("J22").Value = 397.5
$ws.Range("K22").Value = 463.33334
$ws.Range("L22").Value = 397.5
$ws.Range("M22").Value = -290.33334
$ws.Range("N22").Value = -743.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 6326.6665
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 6326.6665
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6326.6665
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -6998.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1132.1562
$ws.Range("I80").Value = 640.7143
$ws.Range("J80").Value = 1514.3889
$ws.Range("K80").Value = 640.7143
$ws.Range("L80").Value = 1514.3889
$ws.Range("M80").Value = 357.2857
$ws.Range("N80").Value = -3510.3889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1132.1562
$ws.Range("I83").Value = 640.7143
$ws.Range("J83").Value = 1514.3889
$ws.Range("K83").Value = 3203.5715
$ws.Range("L83").Value = 7571.9445
$ws.Range("M83").Value = 1788.4285
$ws.Range("N83").Value = -17555.9445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 183851.81
$ws.Range("I105").Value = 202517.8
$ws.Range("J105").Value = 168296.83
$ws.Range("K105").Value = 202517.8
$ws.Range("L105").Value = 168296.83
$ws.Range("M105").Value = -200770.8
$ws.Range("N105").Value = -171790.83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46212.484
$ws.Range("I31").Value = 2458.25
$ws.Range("J31").Value = 69040.78
$ws.Range("K31").Value = 2458.25
$ws.Range("L31").Value = 69040.78
$ws.Range("M31").Value = -2163.25
$ws.Range("N31").Value = -69630.78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 46212.484
$ws.Range("I34").Value = 2458.25
$ws.Range("J34").Value = 69040.78
$ws.Range("K34").Value = 2458.25
$ws.Range("L34").Value = 69040.78
$ws.Range("M34").Value = -2256.25
$ws.Range("N34").Value = -69444.78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1137
$ws.Range("I122").Value = 1066.3334
$ws.Range("J122").Value = 1179.4
$ws.Range("K122").Value = 3199.0002
$ws.Range("L122").Value = 3538.2
$ws.Range("M122").Value = -749.0001999999999
$ws.Range("N122").Value = -8438.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1350
$ws.Range("I34").Value = 100
$ws.Range("J34").Value = 1600
$ws.Range("K34").Value = 300
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -216
$ws.Range("N34").Value = -4968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 460
$ws.Range("I40").Value = 450
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1731
$ws.Range("N40").Value = -2138

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1181.5883
$ws.Range("I113").Value = 1703.1111
$ws.Range("J113").Value = 594.875
$ws.Range("K113").Value = 5109.3333
$ws.Range("L113").Value = 1784.625
$ws.Range("M113").Value = -2939.3333
$ws.Range("N113").Value = -6124.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 667274.3
$ws.Range("J131").Value = 704303.4399999999
$ws.Range("L131").Value = 2112910.32
$ws.Range("N131").Value = -2122990.32

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 49955.375
$ws.Range("I137").Value = 103236
$ws.Range("J137").Value = 11897.786
$ws.Range("K137").Value = 309708
$ws.Range("L137").Value = 35693.358
$ws.Range("M137").Value = -304608
$ws.Range("N137").Value = -45893.358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 66816920
$ws.Range("I80").Value = 111356780
$ws.Range("K80").Value = 111356780
$ws.Range("M80").Value = -111355782

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 66816920
$ws.Range("I83").Value = 111356780
$ws.Range("K83").Value = 556783900
$ws.Range("M83").Value = -556778908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3382.6487
$ws.Range("I132").Value = 2510.2917
$ws.Range("J132").Value = 4993.154
$ws.Range("K132").Value = 7530.875100000001
$ws.Range("L132").Value = 14979.462
$ws.Range("M132").Value = -5000.875100000001
$ws.Range("N132").Value = -20039.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 17199.834
$ws.Range("J38").Value = 17199.834
$ws.Range("L38").Value = 17199.834
$ws.Range("N38").Value = -18019.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1049.4166
$ws.Range("I55").Value = 2094.2856
$ws.Range("J55").Value = 619.17645
$ws.Range("K55").Value = 2094.2856
$ws.Range("L55").Value = 619.17645
$ws.Range("M55").Value = -1921.2856
$ws.Range("N55").Value = -965.17645

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17730.11
$ws.Range("I132").Value = 2202.7173
$ws.Range("J132").Value = 57411.223
$ws.Range("K132").Value = 6608.151899999999
$ws.Range("L132").Value = 172233.669
$ws.Range("M132").Value = -4078.151899999999
$ws.Range("N132").Value = -177293.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3028.0469
$ws.Range("I136").Value = 4141.483
$ws.Range("J136").Value = 2105.4856
$ws.Range("K136").Value = 12424.449
$ws.Range("L136").Value = 6316.4568
$ws.Range("M136").Value = -9874.449000000001
$ws.Range("N136").Value = -11416.4568
